# TC01_CDS_phs001287_Gender_Female.xlsx
# "Added testcases for C3DC phs00468"
#
# The SamplesTab query (cell B3) is rewritten: the two output column
# aliases are renamed ("Tumor" -> "Sample Tumor Status", "Analyte Type" ->
# "Sample Type") and a stray JOIN clause is dropped. Re-typing the cell
# content also nudges Excel into allocating a fresh font/style pair for
# just that cell (it keeps the same visual sz-12 look, but becomes its
# own distinct style entry), so we touch the font explicitly as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newQuery = @'
SELECT
    DISTINCT (smp.sample_id) AS "Sample ID",
    sp.participant_id AS "Participant ID", 
    s.study_name AS "Study Name",
    s.phs_accession AS Accession,
    smp.sample_tumor_status AS "Sample Tumor Status",
    smp.sample_type AS "Sample Type"
FROM 
    df_participant sp
JOIN 
    df_study s ON sp."study.phs_accession" = s.phs_accession
JOIN
    df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_program p ON p.program_acronym = s."program.program_acronym"
JOIN
    df_file f1 ON f1."sample.sample_id" = smp.sample_id
JOIN
    df_genomic_info gi ON gi."file.file_id" = f1.file_id
WHERE 
    s.phs_accession = 'phs001287' AND sp.gender = 'Female'
ORDER BY 
    smp.sample_id ASC
LIMIT 100;
'@

$cell = $ws.Range("B3")
$cell.Value = $newQuery
$cell.WrapText = $true
$cell.Font.ThemeColor = 1
$cell.Font.Size = 12

# User re-selected C3 after editing the SamplesTab query cell.
$ws.Range("C3").Select()
